$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.613.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -7.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.688.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -6.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.58"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -6.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4982"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -16.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2608"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -6.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.65"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -7.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06108"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -10.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07265"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.661.50"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -7.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.435"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -6.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5727"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -8.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.915.54"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008266"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -11.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.75"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -14.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.630.19"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -7.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.005"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -8.80%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.70"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "182.90"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -13.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.182"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -10.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.006"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.89"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -6.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.556"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1132"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -11.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.29"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -7.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.318"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -8.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05617"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -9.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.325"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -6.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.472"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -8.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.466"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -7.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.650"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.006"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.384"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5863"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -8.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.605"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01577"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -7.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.069.64"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.882"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -8.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8484"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.95"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.841.94"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.25"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -7.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000106"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.042"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4328"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05179"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.36%  "
